# Generate Report for Handoff
# Updates the localization-status workbook to reflect that the
# c6e1dc52-95da-4f96-84be-ceed4a84112a.md file is ready for handoff again
# (its handback was stale / not the latest version), across the Overview,
# zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3ccb21ffe8cf6270c48f873e83687bcff3beeae7/e2e/c6e1dc52-95da-4f96-84be-ceed4a84112a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6011cae63d6c8fe087946579be2651b04a3386fd/e2e/c6e1dc52-95da-4f96-84be-ceed4a84112a.md."

# ---- Overview sheet: row 3 is the c6e1dc52-... file ----
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-08-19 08:46:15"

# ---- zh-cn sheet: row 3 is the c6e1dc52-... file ----
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("H3").Value = "2016-08-19 08:46:10"
$zhcn.Range("P3").Value = $errorDetail
$zhcn.Columns.Item(16).ColumnWidth = 39.1

# ---- de-de sheet: row 3 is the c6e1dc52-... file ----
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("H3").Value = "2016-08-19 08:46:15"
$dede.Range("P3").Value = $errorDetail
$dede.Columns.Item(16).ColumnWidth = 39.1
